$d = $word.ActiveDocument

# 1. Remove the _GoBack bookmark from its current location (right after the
#    "The main class is..." paragraph). We'll re-create it at the end of the
#    document, matching its new position in the target layout.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# 2. Append two new paragraphs at the very end of the document body (right
#    before the section break) -- one blank paragraph, then a second blank
#    paragraph that will hold the relocated bookmark.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$tailRange = $lastPara.Range
$tailRange.Collapse(0)
$tailRange.InsertParagraphAfter()
$tailRange.Collapse(0)
$tailRange.InsertParagraphAfter()

# 3. Re-insert the _GoBack bookmark, collapsed, inside the new final paragraph.
$finalPara = $d.Paragraphs($d.Paragraphs.Count)
$finalRange = $finalPara.Range
$finalRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $finalRange)

# 4. Update the sentence describing the fallback behaviour of XMLManager.
$d.Content.Find.Execute(
    "It should also open automatically in your computer, but if not, the error will give you the direction to check the HTML.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "It should also open automatically in your computer, but if it doesn't, the error will give you the direction to check the HTML.",
    2
)

Write-Output "done"
